# Update gh-pages to output generated at 456a3b4
# Applies updated "want to go" (F) and "min ticket price" (G) numbers
# to the 展览 (Exhibitions) sheet and the 全部类型 (All types) sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value = 1783
$ws1.Range("G3").Value = 70

$ws1.Range("F6").Value = 1157

$ws1.Range("F12").Value = 3139

$ws1.Range("F13").Value = 685

$ws1.Range("F19").Value = 1514

$ws1.Range("F20").Value = 308

$ws1.Range("F22").Value = 31

$ws1.Range("F23").Value = 1303

$ws1.Range("F24").Value = 424

$ws1.Range("F25").Value = 518

$ws1.Range("F27").Value = 6850

$ws1.Range("F28").Value = 7162

$ws1.Range("F31").Value = 1720

$ws1.Range("F32").Value = 96

$ws1.Range("F33").Value = 246

# --- Sheet: 全部类型 (All types) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F5").Value = 1783
$ws4.Range("G5").Value = 70

$ws4.Range("F8").Value = 1157

$ws4.Range("F15").Value = 3139

$ws4.Range("F16").Value = 685

$ws4.Range("F22").Value = 1514

$ws4.Range("F23").Value = 308

$ws4.Range("F26").Value = 31

$ws4.Range("F28").Value = 1303

$ws4.Range("F29").Value = 424

$ws4.Range("F30").Value = 518

$ws4.Range("F32").Value = 6850

$ws4.Range("F33").Value = 7162

$ws4.Range("F36").Value = 1720

$ws4.Range("F39").Value = 96

$ws4.Range("F40").Value = 246
